$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuille2")

# Move the "accuracy_Q27" header from J1 to N1.
$ws.Range("N1").Value2 = $ws.Range("J1").Value2
$ws.Range("J1").ClearContents()

# Replicate the accuracy column (J) into the new K, L, M, N columns for
# each data row (2-14).
for ($r = 2; $r -le 14; $r++) {
    $val = $ws.Cells.Item($r, 10).Value2
    $ws.Cells.Item($r, 11).Value2 = $val
    $ws.Cells.Item($r, 12).Value2 = $val
    $ws.Cells.Item($r, 13).Value2 = $val
    $ws.Cells.Item($r, 14).Value2 = $val
}
